# Update "gh-pages output generated at 456a3b4" — refreshed "想去人数"
# (interested-attendee counts) scraped into column F of the 展览 / 演出 /
# 全部类型 sheets. 本地生活 has no changes in this refresh.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 5284
$ws.Cells.Item(4, 6).Value = 5284
$ws.Cells.Item(5, 6).Value = 174
$ws.Cells.Item(6, 6).Value = 217
$ws.Cells.Item(8, 6).Value = 184
$ws.Cells.Item(9, 6).Value = 8833
$ws.Cells.Item(11, 6).Value = 640
$ws.Cells.Item(13, 6).Value = 2595
$ws.Cells.Item(14, 6).Value = 6341
$ws.Cells.Item(15, 6).Value = 2338
$ws.Cells.Item(17, 6).Value = 11
$ws.Cells.Item(18, 6).Value = 28
$ws.Cells.Item(19, 6).Value = 2544
$ws.Cells.Item(22, 6).Value = 6556
$ws.Cells.Item(23, 6).Value = 216
$ws.Cells.Item(24, 6).Value = 81
$ws.Cells.Item(25, 6).Value = 148
$ws.Cells.Item(28, 6).Value = 7148
$ws.Cells.Item(31, 6).Value = 237
$ws.Cells.Item(32, 6).Value = 39
$ws.Cells.Item(36, 6).Value = 22
$ws.Cells.Item(37, 6).Value = 8
$ws.Cells.Item(40, 6).Value = 58
$ws.Cells.Item(41, 6).Value = 2549
$ws.Cells.Item(44, 6).Value = 11
$ws.Cells.Item(47, 6).Value = 553
$ws.Cells.Item(48, 6).Value = 3566
$ws.Cells.Item(50, 6).Value = 1134

# --- 演出 (sheet 2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 20
$ws.Cells.Item(5, 6).Value = 205
$ws.Cells.Item(7, 6).Value = 90
$ws.Cells.Item(15, 6).Value = 160

# --- 全部类型 (sheet 4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 5284
$ws.Cells.Item(4, 6).Value = 5284
$ws.Cells.Item(5, 6).Value = 174
$ws.Cells.Item(6, 6).Value = 218
$ws.Cells.Item(8, 6).Value = 184
$ws.Cells.Item(9, 6).Value = 8833
$ws.Cells.Item(10, 6).Value = 8833
$ws.Cells.Item(12, 6).Value = 640
$ws.Cells.Item(14, 6).Value = 20
$ws.Cells.Item(15, 6).Value = 2595
$ws.Cells.Item(18, 6).Value = 205
$ws.Cells.Item(19, 6).Value = 6341
$ws.Cells.Item(20, 6).Value = 2338
$ws.Cells.Item(21, 6).Value = 90
$ws.Cells.Item(22, 6).Value = 11
$ws.Cells.Item(23, 6).Value = 2544
$ws.Cells.Item(27, 6).Value = 6556
$ws.Cells.Item(28, 6).Value = 216
$ws.Cells.Item(29, 6).Value = 81
$ws.Cells.Item(30, 6).Value = 148
$ws.Cells.Item(32, 6).Value = 7148
$ws.Cells.Item(34, 6).Value = 237
$ws.Cells.Item(36, 6).Value = 22
$ws.Cells.Item(37, 6).Value = 8
$ws.Cells.Item(40, 6).Value = 58
$ws.Cells.Item(41, 6).Value = 2549
$ws.Cells.Item(42, 6).Value = 11
$ws.Cells.Item(45, 6).Value = 553
$ws.Cells.Item(46, 6).Value = 160
$ws.Cells.Item(47, 6).Value = 3566
$ws.Cells.Item(50, 6).Value = 1134
